$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Treatment"
$ws.Range("H1").Value = "Mean NCV (m/s)"
$ws.Range("I1").Value = "Median NCV (m/s)"

$ws.Range("G2").Value = "ctrl"
$ws.Range("H2").Value = 0.34858
$ws.Range("I2").Value = 0.34336

$ws.Range("G3").Value = "ptx200nm"
$ws.Range("H3").Value = 0.2821
$ws.Range("I3").Value = 0.27985

$ws.Range("G1:I3").Select()
